# Generate Report for Handback
# Adds a new handback row (e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2) to the
# Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = 2
$hyperlinkColor = 15570276   # RGB(0x64,0x95,0xED) == FF6495ED, matches existing HyperLink style
$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview": add row 3
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md"
$wsOverview.Range("B3").Value = "e2e\e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7464524c8766c28c6f310edab1ece23869c33d7d/e2e/e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md", "", "", "e2e\e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md") | Out-Null
$wsOverview.Range("B3").Font.Underline = $hyperlinkUnderline
$wsOverview.Range("B3").Font.Color = $hyperlinkColor
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-31 07:26:16"
$wsOverview.Range("G3").NumberFormat = $dateFormat

# ---------------------------------------------------------------------
# Sheet "zh-cn": add row 3
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = "e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7464524c8766c28c6f310edab1ece23869c33d7d/e2e/e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md", "", "", "e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md") | Out-Null
$wsZhCn.Range("A3").Font.Underline = $hyperlinkUnderline
$wsZhCn.Range("A3").Font.Color = $hyperlinkColor
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'True"
$wsZhCn.Range("G3").Value = "e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.2bc5281c51cbe4d2572424672aeda7a73794d73f.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-31 07:25:58"
$wsZhCn.Range("H3").NumberFormat = $dateFormat
$wsZhCn.Range("I3").Value = "e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/af499ece8d1c7414929087b43e61a9cdf1844538/e2e/e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md", "", "", "e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md") | Out-Null
$wsZhCn.Range("I3").Font.Underline = $hyperlinkUnderline
$wsZhCn.Range("I3").Font.Color = $hyperlinkColor
$wsZhCn.Range("J3").Value = "e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.2bc5281c51cbe4d2572424672aeda7a73794d73f.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-31 07:26:50"
$wsZhCn.Range("K3").NumberFormat = $dateFormat
$wsZhCn.Range("L3").Value = "'"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = "'"

# ---------------------------------------------------------------------
# Sheet "de-de": add row 3
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = "e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7464524c8766c28c6f310edab1ece23869c33d7d/e2e/e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md", "", "", "e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md") | Out-Null
$wsDeDe.Range("A3").Font.Underline = $hyperlinkUnderline
$wsDeDe.Range("A3").Font.Color = $hyperlinkColor
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = "e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.2bc5281c51cbe4d2572424672aeda7a73794d73f.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-31 07:26:16"
$wsDeDe.Range("H3").NumberFormat = $dateFormat
$wsDeDe.Range("I3").Value = "e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/05c40c51c315fd9c6b6d0a47ca802f01f6e408af/e2e/e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md", "", "", "e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.md") | Out-Null
$wsDeDe.Range("I3").Font.Underline = $hyperlinkUnderline
$wsDeDe.Range("I3").Font.Color = $hyperlinkColor
$wsDeDe.Range("J3").Value = "e1dcf69c-b114-43f8-9ffa-c2a33a4f25c2.2bc5281c51cbe4d2572424672aeda7a73794d73f.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-31 07:27:14"
$wsDeDe.Range("K3").NumberFormat = $dateFormat
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""
